# Apply cryptos list price/volume update (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.111.35'
$ws.Range("E2").Value = '  -1.26%  '
$ws.Range("D3").Value = '1.855.11'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.80%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4696'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.87%  '
$ws.Range("E8").Value = '  -1.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06521'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.83'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07787'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.49'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.09%  '
$ws.Range("D13").Value = '1.857.56'
$ws.Range("E13").Value = '  -3.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.081'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6608'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '280.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.86%  '
$ws.Range("D17").Value = '30.143.96'
$ws.Range("E17").Value = '  -1.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.449'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.93%  '
$ws.Range("D21").Value = '2.100.84'
$ws.Range("E21").Value = '  -2.74%  '
$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("B23").Value = 'ShibaInu'
$ws.Range("C23").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.000007207'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.108'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '167.75'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.246'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.88'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.906'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -8.37%  '
$ws.Range("E29").Value = '  -3.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09567'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.391'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.466'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.073'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04654'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.093'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6926'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.43%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.716'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01841'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.263'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("E41").Value = '  -4.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '71.39'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8542'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.928'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.07%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.85'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4128'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.83%  '
$ws.Range("D48").Value = '1.010.14'
$ws.Range("E48").Value = '  +4.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.152'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.921'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.96%  '
